$wb = $excel.ActiveWorkbook

# --- 1. Status text update: "Ready for handoff" -> "Handback transform failed" ---
# All cells that previously read "Ready for handoff" shared this one string, so
# every one of them (Overview!E3, Overview!F3, zh-cn!C3, de-de!C3) needs to be
# re-pointed at the new text to keep them in sync, same as an in-place edit of
# the shared string would.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- 2. Populate "Error Detail" (column P) on the zh-cn localization sheet ---
$wsZhCn.Range("P3").Value = "Handback file name: bkgtgobc.t1x is different with handoff file name: e4ea8ccd-75b8-4a1f-9e99-aee63459ac97.34f939f126d1dd695aae1a337b5244d5f8cae21a.zh-cn."
# Widen column P so the error detail text is readable.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- 3. Populate "Error Detail" (column P) on the de-de localization sheet ---
$wsDeDe.Range("P3").Value = "Handback file name: bkgtgobc.t1x is different with handoff file name: e4ea8ccd-75b8-4a1f-9e99-aee63459ac97.34f939f126d1dd695aae1a337b5244d5f8cae21a.de-de."
# Widen column P so the error detail text is readable.
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
